$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 20 with Kenya data (correcting/adding Kenya entry)
$ws.Range("A20").Value = "KEN1"
$ws.Range("B20").Value = "ageband"
$ws.Range("C20").Value = "data/derived/KEN/KEN_agebands.RDS"
$ws.Range("D20").Value = "stratified"
$ws.Range("E20").Value = "aggregate"

# Update selection to match the saved view (cell below the new data)
$ws.Range("A21").Select()
